$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 846.4286
$ws.Range("I32").Value = 550.4
$ws.Range("J32").Value = 1010.8889
$ws.Range("K32").Value = 550.4
$ws.Range("L32").Value = 1010.8889
$ws.Range("M32").Value = -224.4
$ws.Range("N32").Value = -1662.8889
$ws.Range("H125").Value = 565
$ws.Range("I125").Value = 450
$ws.Range("J125").Value = 641.6667
$ws.Range("K125").Value = 4050
$ws.Range("L125").Value = 5775.0003
$ws.Range("M125").Value = -1590
$ws.Range("N125").Value = -10695.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4653.7144
$ws.Range("I46").Value = 5638
$ws.Range("J46").Value = 4260
$ws.Range("K46").Value = 5638
$ws.Range("L46").Value = 4260
$ws.Range("M46").Value = -5319
$ws.Range("N46").Value = -4898
$ws.Range("H61").Value = 1398.2667
$ws.Range("I61").Value = 1459.5385
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 1459.5385
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -1247.5385
$ws.Range("N61").Value = -1424
$ws.Range("H117").Value = 29700
$ws.Range("J117").Value = 29700
$ws.Range("L117").Value = 29700
$ws.Range("N117").Value = -38878
$ws.Range("H122").Value = 23375.682
$ws.Range("I122").Value = 2017.6765
$ws.Range("J122").Value = 79235.08
$ws.Range("K122").Value = 6053.029500000001
$ws.Range("L122").Value = 237705.24
$ws.Range("M122").Value = -3603.029500000001
$ws.Range("N122").Value = -242605.24
$ws.Range("H132").Value = 33368874
$ws.Range("I132").Value = 50001590
$ws.Range("J132").Value = 103441.2
$ws.Range("K132").Value = 150004770
$ws.Range("L132").Value = 310323.6
$ws.Range("M132").Value = -150002240
$ws.Range("N132").Value = -315383.6
$ws.Range("H136").Value = 1398.2667
$ws.Range("I136").Value = 1459.5385
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 4378.6155
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -1828.6155
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7422.483
$ws.Range("I134").Value = 2649.4783
$ws.Range("J134").Value = 25719
$ws.Range("K134").Value = 7948.4349
$ws.Range("L134").Value = 77157
$ws.Range("M134").Value = -5413.4349
$ws.Range("N134").Value = -82227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1690.8334
$ws.Range("I16").Value = 1448
$ws.Range("J16").Value = 1864.2858
$ws.Range("K16").Value = 1448
$ws.Range("L16").Value = 1864.2858
$ws.Range("M16").Value = -1161
$ws.Range("N16").Value = -2438.2858
$ws.Range("H31").Value = 2265627.5
$ws.Range("I31").Value = 3573508.2
$ws.Range("J31").Value = 85826
$ws.Range("K31").Value = 3573508.2
$ws.Range("L31").Value = 85826
$ws.Range("M31").Value = -3573213.2
$ws.Range("N31").Value = -86416
$ws.Range("H34").Value = 2265627.5
$ws.Range("I34").Value = 3573508.2
$ws.Range("J34").Value = 85826
$ws.Range("K34").Value = 3573508.2
$ws.Range("L34").Value = 85826
$ws.Range("M34").Value = -3573306.2
$ws.Range("N34").Value = -86230
$ws.Range("H113").Value = 1690.8334
$ws.Range("I113").Value = 1448
$ws.Range("J113").Value = 1864.2858
$ws.Range("K113").Value = 1448
$ws.Range("L113").Value = 1864.2858
$ws.Range("M113").Value = 722
$ws.Range("N113").Value = -6204.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 799.65515
$ws.Range("I38").Value = 462.9
$ws.Range("J38").Value = 976.8946999999999
$ws.Range("K38").Value = 1388.7
$ws.Range("L38").Value = 2930.6841
$ws.Range("M38").Value = -1041.7
$ws.Range("N38").Value = -3624.6841
$ws.Range("H80").Value = 2740
$ws.Range("J80").Value = 2740
$ws.Range("L80").Value = 8220
$ws.Range("N80").Value = -10092
$ws.Range("H83").Value = 2740
$ws.Range("J83").Value = 2740
$ws.Range("L83").Value = 24660
$ws.Range("N83").Value = -34020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2330.7693
$ws.Range("I97").Value = 1880
$ws.Range("J97").Value = 2612.5
$ws.Range("K97").Value = 1880
$ws.Range("L97").Value = 2612.5
$ws.Range("M97").Value = -1384
$ws.Range("N97").Value = -3604.5
$ws.Range("H132").Value = 36744.69
$ws.Range("I132").Value = 2018.0667
$ws.Range("J132").Value = 73951.78999999999
$ws.Range("K132").Value = 6054.2001
$ws.Range("L132").Value = 221855.37
$ws.Range("M132").Value = -3524.2001
$ws.Range("N132").Value = -226915.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1875.3572
$ws.Range("I7").Value = 1896.25
$ws.Range("J7").Value = 1750
$ws.Range("K7").Value = 1896.25
$ws.Range("L7").Value = 1750
$ws.Range("M7").Value = -1784.25
$ws.Range("N7").Value = -1974
$ws.Range("H40").Value = 1837.2084
$ws.Range("I40").Value = 1729
$ws.Range("J40").Value = 2100
$ws.Range("K40").Value = 1729
$ws.Range("L40").Value = 2100
$ws.Range("M40").Value = -1593
$ws.Range("N40").Value = -2372
$ws.Range("H126").Value = 1875.3572
$ws.Range("I126").Value = 1896.25
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 5688.75
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -3218.75
$ws.Range("N126").Value = -10190
$ws.Range("H132").Value = 2943620
$ws.Range("I132").Value = 5557402
$ws.Range("J132").Value = 3115.1875
$ws.Range("K132").Value = 16672206
$ws.Range("L132").Value = 9345.5625
$ws.Range("M132").Value = -16669676
$ws.Range("N132").Value = -14405.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18292.334
$ws.Range("I41").Value = 30000
$ws.Range("J41").Value = 12438.5
$ws.Range("K41").Value = 30000
$ws.Range("L41").Value = 12438.5
$ws.Range("M41").Value = -29610
$ws.Range("N41").Value = -13218.5
$ws.Range("H45").Value = 10846.818
$ws.Range("J45").Value = 10846.818
$ws.Range("L45").Value = 10846.818
$ws.Range("N45").Value = -11828.818
$ws.Range("H123").Value = 34354.92
$ws.Range("J123").Value = 34354.92
$ws.Range("L123").Value = 34354.92
$ws.Range("N123").Value = -44154.92
